$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Odd_Over05_HT (P) for row 2
$ws.Range("P2").Value = 1.47

# Update Odd_Under25_FT (O), Odd_Over05_HT (P), Odd_Under05_HT (Q) for row 3
$ws.Range("O3").Value = 1.58
$ws.Range("P3").Value = 1.41
$ws.Range("Q3").Value = 2.62
